$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: update values to the new "custom accuracy" (rounded) figures
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH")
$vals = @(10.09,7.23,0.9,21.63,17.88,7.94,34.08,12.22,5.34,7.97,8.779999999999999,9.140000000000001,2.54,7.9,11.18,6.78,0.78,0.53,112.86,22.2,7.29,14.77,7.91,1.14,15.93,6.44,5.82,6.8,9.15,0.5600000000000001,30.99,4.04,9.109999999999999)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $vals[$i]
}

# Row 6 is no longer part of the (now 1000-row-limited) dataset - remove it entirely
$ws.Rows.Item(6).Delete()
